$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data from GitHub Actions scrape

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.044.80'
$ws.Range("E2").Value = '  -3.32%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.600.22'
$ws.Range("E3").Value = '  -2.15%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.01%  '

# Row 5
$ws.Range("E5").Value = '  -0.02%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '301.77'
$ws.Range("E6").Value = '  -2.16%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3780'
$ws.Range("E7").Value = '  -1.72%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3645'
$ws.Range("E8").Value = '  -4.04%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '49.88'
$ws.Range("E9").Value = '  -1.55%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.260'
$ws.Range("E10").Value = '  -4.61%  '

# Row 11 (Dogecoin)
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08140'
$ws.Range("E11").Value = '  -2.55%  '

# Row 12 (BinanceUSD)
$ws.Range("B12").Value = 'BinanceUSD'
$ws.Range("C12").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.000'
$ws.Range("E12").Value = '  +0.02%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.55'
$ws.Range("E13").Value = '  -4.54%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.590'
$ws.Range("E14").Value = '  -5.00%  '

# Row 15 (ShibaInu)
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001248'
$ws.Range("E15").Value = '  -3.89%  '

# Row 16 (Chainlink)
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.350'
$ws.Range("E16").Value = '  -5.02%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.606.96'
$ws.Range("E17").Value = '  -1.62%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '91.97'
$ws.Range("E18").Value = '  -1.44%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06826'
$ws.Range("E19").Value = '  -1.39%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.24'
$ws.Range("E20").Value = '  -5.64%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.538'
$ws.Range("E21").Value = '  -4.27%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.5578'
$ws.Range("E22").Value = '  -5.04%  '

# Row 23
$ws.Range("E23").Value = '  +0.13%  '

# Row 24
$ws.Range("E24").Value = '  -3.10%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '23.040.22'
$ws.Range("E25").Value = '  -3.31%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.366'
$ws.Range("E26").Value = '  -2.44%  '

# Row 27
$ws.Range("E27").Value = '  -1.07%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.04'
$ws.Range("E28").Value = '  -3.21%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '150.43'
$ws.Range("E29").Value = '  -1.25%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.232'
$ws.Range("E30").Value = '  -3.98%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '134.20'
$ws.Range("E31").Value = '  -1.37%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.355'
$ws.Range("E32").Value = '  -4.85%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.823'
$ws.Range("E33").Value = '  -12.92%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.782.83'
$ws.Range("E34").Value = '  -2.39%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9649'
$ws.Range("E35").Value = '  -1.37%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.07582'
$ws.Range("E36").Value = '  -4.08%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '10.35'
$ws.Range("E37").Value = '  -0.40%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.259'
$ws.Range("E38").Value = '  -4.76%  '

# Row 39
$ws.Range("E39").Value = '  -5.84%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2527'
$ws.Range("E40").Value = '  -4.44%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.08874'
$ws.Range("E41").Value = '  -1.95%  '

# Row 42
$ws.Range("E42").Value = '  -3.29%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7032'
$ws.Range("E43").Value = '  -5.62%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '12.39'
$ws.Range("E44").Value = '  -6.06%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '15.25'
$ws.Range("E45").Value = '  -7.89%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6628'
$ws.Range("E46").Value = '  -3.26%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.292'
$ws.Range("E48").Value = '  -4.47%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.993'
$ws.Range("E49").Value = '  -1.67%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '132.58'
$ws.Range("E50").Value = '  -1.15%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07902'
$ws.Range("E51").Value = '  -3.57%  '
